$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1417.875
$ws.Range("I19").Value = 891.3333
$ws.Range("K19").Value = 891.3333
$ws.Range("M19").Value = -716.3333
$ws.Range("H86").Value = 7827
$ws.Range("I86").Value = 7699
$ws.Range("J86").Value = 7891
$ws.Range("K86").Value = 7699
$ws.Range("L86").Value = 7891
$ws.Range("M86").Value = -6576
$ws.Range("N86").Value = -10137
$ws.Range("H89").Value = 7827
$ws.Range("I89").Value = 7699
$ws.Range("J89").Value = 7891
$ws.Range("K89").Value = 38495
$ws.Range("L89").Value = 39455
$ws.Range("M89").Value = -32879
$ws.Range("N89").Value = -50687
$ws.Range("H103").Value = 2541.3333
$ws.Range("J103").Value = 3418.25
$ws.Range("L103").Value = 10254.75
$ws.Range("N103").Value = -11426.75
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("H131").Value = 2580
$ws.Range("I131").Value = 2580
$ws.Range("K131").Value = 7740
$ws.Range("M131").Value = -2700
$ws.Range("H137").Value = 2766.6667
$ws.Range("I137").Value = 2400
$ws.Range("K137").Value = 7200
$ws.Range("M137").Value = -4650

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3523.1738
$ws.Range("I32").Value = 2906.3333
$ws.Range("K32").Value = 2906.3333
$ws.Range("M32").Value = -2619.3333
$ws.Range("H74").Value = 993.7
$ws.Range("I74").Value = 1017.125
$ws.Range("K74").Value = 1017.125
$ws.Range("M74").Value = -143.125
$ws.Range("H77").Value = 993.7
$ws.Range("I77").Value = 1017.125
$ws.Range("K77").Value = 5085.625
$ws.Range("M77").Value = -717.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 468.66666
$ws.Range("J64").Value = 419
$ws.Range("L64").Value = 419
$ws.Range("N64").Value = -869
$ws.Range("H67").Value = 468.66666
$ws.Range("J67").Value = 419
$ws.Range("L67").Value = 419
$ws.Range("N67").Value = -1979
$ws.Range("H86").Value = 3712.75
$ws.Range("I86").Value = 3728.111
$ws.Range("J86").Value = 3666.6667
$ws.Range("K86").Value = 3728.111
$ws.Range("L86").Value = 3666.6667
$ws.Range("M86").Value = -2605.111
$ws.Range("N86").Value = -5912.6667
$ws.Range("H89").Value = 3712.75
$ws.Range("I89").Value = 3728.111
$ws.Range("J89").Value = 3666.6667
$ws.Range("K89").Value = 18640.555
$ws.Range("L89").Value = 18333.3335
$ws.Range("M89").Value = -13024.555
$ws.Range("N89").Value = -29565.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5126.067
$ws.Range("I22").Value = 6464.1
$ws.Range("J22").Value = 2450
$ws.Range("K22").Value = 6464.1
$ws.Range("L22").Value = 2450
$ws.Range("M22").Value = -6114.1
$ws.Range("N22").Value = -3150
$ws.Range("H99").Value = 3158.2307
$ws.Range("I99").Value = 2163.8572
$ws.Range("K99").Value = 2163.8572
$ws.Range("M99").Value = -665.8571999999999
$ws.Range("H126").Value = 3158.2307
$ws.Range("I126").Value = 2163.8572
$ws.Range("K126").Value = 6491.571599999999
$ws.Range("M126").Value = -4021.571599999999
$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -85060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 397
$ws.Range("J38").Value = 295.83334
$ws.Range("L38").Value = 887.5000200000001
$ws.Range("N38").Value = -1581.50002
$ws.Range("H69").Value = 1200
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 1200
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H80").Value = 14995
$ws.Range("I80").Value = 14995
$ws.Range("K80").Value = 44985
$ws.Range("M80").Value = -44049
$ws.Range("H83").Value = 14995
$ws.Range("I83").Value = 14995
$ws.Range("K83").Value = 134955
$ws.Range("M83").Value = -130275
$ws.Range("H92").Value = 553
$ws.Range("I92").Value = 441.25
$ws.Range("K92").Value = 1323.75
$ws.Range("M92").Value = -75.75
$ws.Range("H107").Value = 1765.6666
$ws.Range("I107").Value = 1803
$ws.Range("K107").Value = 5409
$ws.Range("M107").Value = -3489
$ws.Range("H114").Value = 458
$ws.Range("I114").Value = 524.5
$ws.Range("K114").Value = 1573.5
$ws.Range("M114").Value = 1680.5
$ws.Range("H117").Value = 1519
$ws.Range("I117").Value = 562.25
$ws.Range("K117").Value = 1686.75
$ws.Range("M117").Value = 1755.25
$ws.Range("H121").Value = 1083.6
$ws.Range("I121").Value = 430
$ws.Range("J121").Value = 1519.3334
$ws.Range("K121").Value = 1290
$ws.Range("L121").Value = 4558.0002
$ws.Range("M121").Value = 20
$ws.Range("N121").Value = -7178.0002
$ws.Range("H131").Value = 977.15625
$ws.Range("J131").Value = 986.129
$ws.Range("L131").Value = 2958.387
$ws.Range("N131").Value = -13038.387

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 196.1579
$ws.Range("I2").Value = 53
$ws.Range("K2").Value = 53
$ws.Range("M2").Value = 60
$ws.Range("H102").Value = 1365.25
$ws.Range("I102").Value = 1182.3334
$ws.Range("K102").Value = 1182.3334
$ws.Range("M102").Value = 439.6666
$ws.Range("H122").Value = 1972.2858
$ws.Range("I122").Value = 2381.4
$ws.Range("J122").Value = 949.5
$ws.Range("K122").Value = 7144.200000000001
$ws.Range("L122").Value = 2848.5
$ws.Range("M122").Value = -4694.200000000001
$ws.Range("N122").Value = -7748.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4966.6665
$ws.Range("I132").Value = 4950
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 14850
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -12320
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 6243.8184
$ws.Range("I136").Value = 6147.125
$ws.Range("J136").Value = 6501.6665
$ws.Range("K136").Value = 18441.375
$ws.Range("L136").Value = 19504.9995
$ws.Range("M136").Value = -15891.375
